$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a2 = @'
INVEMAR - Calidad Ambiental Marina
'@
$a3 = @'
Grupo de Investigación en Modelación y Evaluación de Sistemas Ambientales. GiMESA
'@
$a4 = @'
GRUPO DE INVESTIGACION AMBIENTAL 
'@
$b2 = @'
5.- Extensión extracurricular : DISEÑO DE MONITOREO DE CALIDAD AMBIENTAL EN ZONAS MARINO-COSTERAS: Interpretación de parámetros indicadores de calidad y contaminación de aguas marinas y costeras  Colombia, 2015, Idioma: Español, Medio de divulgación: CD-ROM/DVD Sitio web: , Participación como Docente, Duración (semanas): 1, Finalidad: Brindar conceptos sobre la interpretación de los parámetros indicadores de calidad y contaminación de aguas marinas y costeras, y aplicaciones en casos prácticos que faciliten la gestión ambiental de las CAR costeras, el MADS y otras entidades relacionadas con el manejo del recurso hídrico marino y Lugar: INVEMAR, Institución financiadora: INSTITUTO DE INVESTIGACIONES MARINAS Y COSTERAS JOSE BENITO VIVES DE ANDREIS INVEMAR  Autores: LIZBETH JANET VIVAS AGUAS
'@
$b3 = @'
2.- Extensión extracurricular : MONITOREO AMBIENTAL PARTICIPATIVOO  Colombia, 2021, Idioma: Español, Medio de divulgación: Papel Sitio web: , Participación como Docente, Duración (semanas): 0, Finalidad: Capacitar a las comunidades del área de influencia de CERRO MATOSO en Monitoreos de los componentes agua, suelo, aguas subterráneas , SIG, flora, fauna y aire Lugar: Municipios de Montelíbano, Puerto Libertador y San José de Uré, Institución financiadora: CERROMATOSO S.A.  Autores: VIVIANA CECILIA SOTO BARRERA, GABRIEL ANTONIO CAMPO DAZA, ZORAYA YASQUINE MARTINEZ LARA, DORIS MEJIA AVILA, ANGELICA BUSTAMANTE RUIZ, MONICA CECILIA CANTERO BENITEZ 
 3.- Extensión extracurricular : DIPLOMADO EN MONITOREOAMBIENTAL PARTICIPATIVO  Colombia, 2021, Idioma: Español, Medio de divulgación: Papel Sitio web: , Participación como Docente, Duración (semanas): 12, Finalidad: Capacitar a las comunidades del área de influencia de CERRO MATOSO en Monitoreos de los componentes agua, suelo, aguas subterráneas , SIG, flora, fauna y aire Lugar: MONTELÍBANO, Institución financiadora: CERROMATOSO S.A.  Autores: MONICA CECILIA CANTERO BENITEZ 
 4.- Extensión extracurricular : Capacitar a las comunidades del área de influencia de CERRO MATOSO en Monitoreos de los componentes agua, suelo, aguas subterráneas , SIG, flora, fauna y aire  Colombia, 2020, Idioma: Español, Medio de divulgación: Varios Sitio web: , Participación como Docente, Duración (semanas): 0, Finalidad: Lugar: , Institución financiadora: Autores: ANGELICA BUSTAMANTE RUIZ
'@
$b4 = @'
1.- Extensión extracurricular : Curso corto: Fortalecimiento de prestadores del servicio de acueducto del sector rural  Colombia, 2018, Idioma: Español, Medio de divulgación: Papel Sitio web: , Participación como Docente, Duración (semanas): 3, Finalidad: Capacitar a los integrantes de las juntas administradoras de acueducto en temáticas relacionadas con: Gestión ambiental, Mantenimiento y operación de los componentes del acueducto, monitoreo de la calidad del agua, gestión administrativa y financiera Lugar: Ipiales - Nariño, Institución financiadora: Fondo Rotatorio del Ministerio de Relaciones Exteriores  Autores: PAOLA ANDREA ORTEGA GUERRERO 
 2.- Extensión extracurricular : Curso corto: Fortalecimiento de prestadores del servicio de acueducto del sector rural  Colombia, 2018, Idioma: Español, Medio de divulgación: Papel Sitio web: , Participación como Organizador, Duración (semanas): 3, Finalidad: Capacitar a los integrantes de las juntas administradoras de acueducto en temáticas relacionadas con: Gestión ambiental, Mantenimiento y operación de los componentes del acueducto, monitoreo de la calidad del agua, gestión administrativa y financiera Lugar: Ipiales - Nariño, Institución financiadora: Fondo Rotatorio del Ministerio de Relaciones Exteriores  Autores: PAOLA ANDREA ORTEGA GUERRERO
'@

$ws.Range("A2").Value2 = $a2
$ws.Range("B2").Value2 = $b2
$ws.Range("A3").Value2 = $a3
$ws.Range("B3").Value2 = $b3
$ws.Range("A4").Value2 = $a4
$ws.Range("B4").Value2 = $b4

$ws.Rows(2).EntireRow.AutoFit()
$ws.Rows(3).EntireRow.AutoFit()
$ws.Rows(4).EntireRow.AutoFit()

$ws.Rows(5).Delete()
